# corretto il meccanismo che genera Terminologia_glossario
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = 14000076

$ws.Range("C34").Value = 14000052

$ws.Range("D36").Value = 43844
$ws.Range("D37").Value = 43804
$ws.Range("D38").Value = 43774
$ws.Range("B39").Value = 13000038
$ws.Range("D39").Value = 43845

$ws.Range("B40").Value = 13000039
$ws.Range("C40").Value = 14000039
$ws.Range("B41").Value = 13000039
$ws.Range("C41").Value = 14000040
$ws.Range("B42").Value = 13000039
$ws.Range("C42").Value = 14000041
$ws.Range("B43").Value = 13000039
$ws.Range("C43").Value = 14000043

$ws.Range("B44").Value = 13000040
$ws.Range("C44").Value = 14000057
$ws.Range("D44").Value = 43767
$ws.Range("E44").Value = "ITCH00016"

$ws.Range("B45").Value = 13000040
$ws.Range("C45").Value = 14000058
$ws.Range("D45").Value = 43767
$ws.Range("E45").Value = "ITCH00017"

$ws.Range("B46").Value = 13000040
$ws.Range("C46").Value = 14000059
$ws.Range("D46").Value = 43767
$ws.Range("E46").Value = "ITCH00018"

$ws.Range("B47").Value = 13000040
$ws.Range("C47").Value = 14000060
$ws.Range("D47").Value = 43767
$ws.Range("E47").Value = "ITCH00019"

$ws.Range("B48").Value = 13000040
$ws.Range("C48").Value = 14000061
$ws.Range("D48").Value = 43767
$ws.Range("E48").Value = "ITCH00020"

$ws.Range("B49").Value = 13000040
$ws.Range("C49").Value = 14000063
$ws.Range("D49").Value = 43767
$ws.Range("E49").Value = "ITCH00022"

$ws.Range("B50").Value = 13000040
$ws.Range("C50").Value = 14000064
$ws.Range("D50").Value = 43767
$ws.Range("E50").Value = "ITCH00023"

$ws.Range("B51").Value = 13000040
$ws.Range("C51").Value = 14000065
$ws.Range("D51").Value = 43767
$ws.Range("E51").Value = "ITCH00024"

$ws.Range("B52").Value = 13000040
$ws.Range("C52").Value = 14000071
$ws.Range("D52").Value = 43767
$ws.Range("E52").Value = "ITCH00028"

$ws.Range("B53").Value = 13000044
$ws.Range("C53").Value = 14000044
$ws.Range("D53").Value = 43657
$ws.Range("E53").Value = "ITCH00005"

$ws.Range("B54").Value = 13000044
$ws.Range("C54").Value = 14000045
$ws.Range("D54").Value = 43657
$ws.Range("E54").Value = "ITCH00006"

$ws.Range("B55").Value = 13000046
$ws.Range("C55").Value = 14000046
$ws.Range("D55").Value = 43731
$ws.Range("E55").Value = "ITCH00007"

$ws.Range("B56").Value = 13000047
$ws.Range("C56").Value = 14000047
$ws.Range("D56").Value = 43766
$ws.Range("E56").Value = "ITCH00009"

$ws.Range("B57").Value = 13000047
$ws.Range("C57").Value = 14000049
$ws.Range("D57").Value = 43766
$ws.Range("E57").Value = "ITCH00010"

$ws.Range("B58").Value = 13000047
$ws.Range("C58").Value = 14000051
$ws.Range("D58").Value = 43766
$ws.Range("E58").Value = "ITCH00011"

$ws.Range("B59").Value = 13000047
$ws.Range("C59").Value = 14000054
$ws.Range("D59").Value = 43766
$ws.Range("E59").Value = "ITCH00013"

$ws.Range("B60").Value = 13000047
$ws.Range("C60").Value = 14000055
$ws.Range("D60").Value = 43766
$ws.Range("E60").Value = "ITCH00014"

$ws.Range("B61").Value = 13000047
$ws.Range("C61").Value = 14000056
$ws.Range("D61").Value = 43766
$ws.Range("E61").Value = "ITCH00015"

$ws.Range("B62").Value = 13000047
$ws.Range("C62").Value = 14000062
$ws.Range("D62").Value = 43766
$ws.Range("E62").Value = "ITCH00021"

$ws.Range("B63").Value = 13000047
$ws.Range("C63").Value = 14000074
$ws.Range("D63").Value = 43766
$ws.Range("E63").Value = "ITCH00030"

$ws.Range("B64").Value = 13000066
$ws.Range("C64").Value = 14000066

$ws.Range("B65").Value = 13000068
$ws.Range("C65").Value = 14000068

$ws.Range("B66").Value = 13000068
$ws.Range("C66").Value = 14000073

$ws.Range("B67").Value = 13000068
$ws.Range("C67").Value = 14000075

$ws.Range("B68").Value = 13000069
$ws.Range("C68").Value = 14000069
